$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (header is row 1, data starts row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {

    # 1. Update the "Förändrad" (changed) date in column C from 45184 to 45186
    $curC = $ws.Cells.Item($r, 3).Value2
    if ($curC -eq 45184) {
        $ws.Cells.Item($r, 3).Value2 = 45186
    }

    # Beteckning (case/report id) lives in column A and is used as the
    # friendly display text for the HYPERLINK formulas on this row.
    $beteckning = $ws.Cells.Item($r, 1).Value2
    if ([string]::IsNullOrEmpty($beteckning)) {
        continue
    }

    # 2. Add a second argument (friendly name) to every HYPERLINK() formula
    #    in this row that doesn't already have one, for columns S..Y.
    for ($c = 19; $c -le 25; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f -like 'HYPERLINK(*' -or $f -like '=HYPERLINK(*') {
                if ($f -notlike '*,*') {
                    $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $beteckning + '")'
                    $cell.Formula = $newFormula
                }
            }
        }
    }
}
